$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'61.984.73"
$ws.Range("E2").Value = "  +2.50%  "
$ws.Range("D3").Value = "'2.405.05"
$ws.Range("E3").Value = "  +3.40%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'561.31"
$ws.Range("E5").Value = "  +2.84%  "
$ws.Range("D6").Value = "'138.73"
$ws.Range("E6").Value = "  +5.80%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  +1.13%  "
$ws.Range("D9").Value = "'2.402.19"
$ws.Range("E9").Value = "  +3.35%  "
$ws.Range("E10").Value = "  +3.16%  "
$ws.Range("E11").Value = "  +4.00%  "
$ws.Range("E12").Value = "  +0.09%  "
$ws.Range("E13").Value = "  +4.01%  "
$ws.Range("D14").Value = "'25.64"
$ws.Range("E14").Value = "  +8.15%  "
$ws.Range("D15").Value = "'2.833.10"
$ws.Range("E15").Value = "  +3.36%  "
$ws.Range("D16").Value = "'61.914.94"
$ws.Range("E16").Value = "  +2.43%  "
$ws.Range("D17").Value = "'0.0000139"
$ws.Range("E17").Value = "  +4.22%  "
$ws.Range("D18").Value = "'2.407.87"
$ws.Range("E18").Value = "  +3.74%  "
$ws.Range("D19").Value = "'11.00"
$ws.Range("E19").Value = "  +3.68%  "
$ws.Range("D20").Value = "'341.77"
$ws.Range("E20").Value = "  +8.44%  "
$ws.Range("E21").Value = "  +1.72%  "
$ws.Range("D22").Value = "'6.88"
$ws.Range("E22").Value = "  +3.86%  "
$ws.Range("E23").Value = "  +0.30%  "
$ws.Range("D24").Value = "'64.85"
$ws.Range("E24").Value = "  +1.37%  "
$ws.Range("E25").Value = "  +1.05%  "
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("D27").Value = "'8.33"
$ws.Range("E27").Value = "  +5.96%  "
$ws.Range("E28").Value = "  +11.12%  "
$ws.Range("E29").Value = "  +16.17%  "
$ws.Range("E30").Value = "  +3.95%  "
$ws.Range("D31").Value = "'0.0₃0769"
$ws.Range("E31").Value = "  +4.96%  "
$ws.Range("E32").Value = "  +7.07%  "
$ws.Range("D33").Value = "'171.08"
$ws.Range("E33").Value = "  -1.49%  "
$ws.Range("E34").Value = "  +3.15%  "
$ws.Range("D35").Value = "'1.40"
$ws.Range("E35").Value = "  +1.13%  "
$ws.Range("D36").Value = "'4.56"
$ws.Range("E36").Value = "  +12.51%  "
$ws.Range("D37").Value = "'18.47"
$ws.Range("E37").Value = "  +3.49%  "
$ws.Range("D38").Value = "'363.35"
$ws.Range("E38").Value = "  +11.28%  "
$ws.Range("E40").Value = "  -0.11%  "
$ws.Range("E41").Value = "  +8.02%  "
$ws.Range("D42").Value = "'38.95"
$ws.Range("E42").Value = "  +2.68%  "
$ws.Range("D43").Value = "'143.87"
$ws.Range("E43").Value = "  +4.16%  "
$ws.Range("E44").Value = "  +4.67%  "
$ws.Range("D45").Value = "'20.39"
$ws.Range("E45").Value = "  +5.82%  "
$ws.Range("D46").Value = "'0.0959"
$ws.Range("E46").Value = "  +1.68%  "
$ws.Range("E47").Value = "  +4.04%  "
$ws.Range("D48").Value = "'0.584"
$ws.Range("E48").Value = "  +4.18%  "
$ws.Range("E49").Value = "  +3.99%  "
$ws.Range("D50").Value = "'17.77"
$ws.Range("E50").Value = "  +5.22%  "
$ws.Range("D51").Value = "'0.0₆0217"
$ws.Range("E51").Value = "  +0.64%  "
